$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"24.76000000000043"
$ws.Range("H2").Value = [double]"4.312217249946571e-12"
$ws.Range("I2").Value = [double]"4.312217249946571e-12"
$ws.Range("L2").Value = [double]"60.26446669594556"
$ws.Range("M2").Value = "[46.65140382809224, 73.87752956379887]"
$ws.Range("N2").Value = [double]"1.66378022470326e-11"
$ws.Range("O2").Value = [double]"1.66378022470326e-11"
$ws.Range("P2").Value = [double]"1.616395018964117"
$ws.Range("Q2").Value = "[1.364816027685655, 1.8679740102425795]"
$ws.Range("R2").Value = [double]"2.220446049250313e-16"
$ws.Range("S2").Value = [double]"2.220446049250313e-16"
$ws.Range("T2").Value = [double]"53.74910522302656"
$ws.Range("U2").Value = "[44.9888845409387, 62.50932590511442]"
$ws.Range("V2").Value = [double]"4.440892098500626e-16"
$ws.Range("W2").Value = [double]"4.440892098500626e-16"
$ws.Range("X2").Value = [double]"18.39031031031063"
$ws.Range("Y2").Value = [double]"17.39891891891922"
$ws.Range("Z2").Value = [double]"19.38170170170204"

# Row 3
$ws.Range("F3").Value = [double]"24.76000000000043"
$ws.Range("H3").Value = [double]"6.217248937900877e-15"
$ws.Range("I3").Value = [double]"6.217248937900877e-15"
$ws.Range("L3").Value = [double]"58.39243796071343"
$ws.Range("M3").Value = "[46.086541846239285, 70.69833407518757]"
$ws.Range("N3").Value = [double]"2.108091479158247e-12"
$ws.Range("O3").Value = [double]"2.108091479158247e-12"
$ws.Range("P3").Value = [double]"1.289342330302117"
$ws.Range("Q3").Value = "[1.0629212381515023, 1.5157634224527312]"
$ws.Range("R3").Value = [double]"5.995204332975845e-15"
$ws.Range("S3").Value = [double]"5.995204332975845e-15"
$ws.Range("T3").Value = [double]"54.120657400825"
$ws.Range("U3").Value = "[47.205736095946946, 61.035578705703045]"
$ws.Range("X3").Value = [double]"19.67911911911946"
$ws.Range("Y3").Value = [double]"18.7868668668672"
$ws.Range("Z3").Value = [double]"20.57137137137173"

# Row 4
$ws.Range("F4").Value = [double]"24.76000000000043"
$ws.Range("H4").Value = [double]"2.583488978302739e-13"
$ws.Range("I4").Value = [double]"2.583488978302739e-13"
$ws.Range("L4").Value = [double]"61.913232321692"
$ws.Range("M4").Value = "[47.66963755438782, 76.15682708899618]"
$ws.Range("N4").Value = [double]"2.82167622600582e-11"
$ws.Range("O4").Value = [double]"2.82167622600582e-11"
$ws.Range("P4").Value = [double]"1.301921279866041"
$ws.Range("Q4").Value = "[1.0503422885875784, 1.553500271144503]"
$ws.Range("R4").Value = [double]"1.396660564978447e-13"
$ws.Range("S4").Value = [double]"1.396660564978447e-13"
$ws.Range("T4").Value = [double]"52.54860027007838"
$ws.Range("U4").Value = "[44.35061201124144, 60.746588528915325]"
$ws.Range("V4").Value = [double]"2.220446049250313e-16"
$ws.Range("W4").Value = [double]"2.220446049250313e-16"
$ws.Range("X4").Value = [double]"19.62954954954989"
$ws.Range("Y4").Value = [double]"18.63815815815848"
$ws.Range("Z4").Value = [double]"20.6209409409413"

# Row 5
$ws.Range("F5").Value = [double]"24.76000000000043"
$ws.Range("H5").Value = [double]"1.061369880872576e-10"
$ws.Range("I5").Value = [double]"1.061369880872576e-10"
$ws.Range("L5").Value = [double]"58.90899099554235"
$ws.Range("M5").Value = "[40.93286384702785, 76.88511814405683]"
$ws.Range("N5").Value = [double]"3.997743314165803e-08"
$ws.Range("O5").Value = [double]"3.997743314165803e-08"
$ws.Range("P5").Value = [double]"0.9622896416401163"
$ws.Range("Q5").Value = "[0.6603948521059619, 1.2641844311742707]"
$ws.Range("R5").Value = [double]"7.416098868340271e-08"
$ws.Range("S5").Value = [double]"7.416098868340271e-08"
$ws.Range("T5").Value = [double]"56.94136387705253"
$ws.Range("U5").Value = "[47.55478097803929, 66.32794677606577]"
$ws.Range("V5").Value = [double]"6.661338147750939e-16"
$ws.Range("W5").Value = [double]"6.661338147750939e-16"
$ws.Range("X5").Value = [double]"20.96792792792829"
$ws.Range("Y5").Value = [double]"19.7782582582586"
$ws.Range("Z5").Value = [double]"22.15759759759798"

# Row 6
$ws.Range("F6").Value = [double]"25.78000000000059"
$ws.Range("H6").Value = [double]"4.837574785199195e-12"
$ws.Range("I6").Value = [double]"4.837574785199195e-12"
$ws.Range("L6").Value = [double]"55.25573442890708"
$ws.Range("M6").Value = "[40.152743394949375, 70.35872546286478]"
$ws.Range("N6").Value = [double]"2.901709539315789e-09"
$ws.Range("O6").Value = [double]"2.901709539315789e-09"
$ws.Range("P6").Value = [double]"0.748447499053424"
$ws.Range("Q6").Value = "[0.4842895582110387, 1.0126054398958093]"
$ws.Range("R6").Value = [double]"8.494301764194745e-07"
$ws.Range("S6").Value = [double]"8.494301764194745e-07"
$ws.Range("T6").Value = [double]"53.77337112878034"
$ws.Range("U6").Value = "[46.009154790192085, 61.537587467368596]"
$ws.Range("V6").Value = [double]"0"
$ws.Range("W6").Value = [double]"0"
$ws.Range("X6").Value = [double]"22.70910910910963"
$ws.Range("Y6").Value = [double]"21.62526526526576"
$ws.Range("Z6").Value = [double]"23.7929529529535"

# Row 7
$ws.Range("F7").Value = [double]"25.78000000000059"
$ws.Range("H7").Value = [double]"4.907185768843192e-14"
$ws.Range("I7").Value = [double]"4.907185768843192e-14"
$ws.Range("L7").Value = [double]"64.0292610328674"
$ws.Range("M7").Value = "[50.2086713060349, 77.84985075969989]"
$ws.Range("N7").Value = [double]"4.344524739963163e-12"
$ws.Range("O7").Value = [double]"4.344524739963163e-12"
$ws.Range("P7").Value = [double]"0.3333421634439615"
$ws.Range("Q7").Value = "[0.10692107129334527, 0.5597632555945777]"
$ws.Range("R7").Value = [double]"0.004825303314059726"
$ws.Range("S7").Value = [double]"0.004825303314059726"
$ws.Range("T7").Value = [double]"56.35304360778357"
$ws.Range("U7").Value = "[48.67298983912369, 64.03309737644345]"
$ws.Range("X7").Value = [double]"24.41229229229285"
$ws.Range("Y7").Value = [double]"23.48328328328382"
$ws.Range("Z7").Value = [double]"25.34130130130189"

# Row 8
$ws.Range("F8").Value = [double]"25.78000000000059"
$ws.Range("H8").Value = [double]"1.389999226830696e-13"
$ws.Range("I8").Value = [double]"1.389999226830696e-13"
$ws.Range("L8").Value = [double]"63.80360675370871"
$ws.Range("M8").Value = "[51.24117461331906, 76.36603889409837]"
$ws.Range("N8").Value = [double]"2.542410726391608e-13"
$ws.Range("O8").Value = [double]"2.542410726391608e-13"
$ws.Range("P8").Value = [double]"-0.0503157982556921"
$ws.Range("Q8").Value = "[-0.25786846606042246, 0.15723686954903826]"
$ws.Range("R8").Value = [double]"0.6277312811237863"
$ws.Range("S8").Value = [double]"0.6277312811237863"
$ws.Range("T8").Value = [double]"61.76932896621931"
$ws.Range("U8").Value = "[53.80341070927143, 69.73524722316718]"
$ws.Range("X8").Value = [double]"0.2064464464464528"
$ws.Range("Y8").Value = [double]"-0.6451451451451554"
$ws.Range("Z8").Value = [double]"1.058038038038061"

# Row 9
$ws.Range("F9").Value = [double]"25.78000000000059"
$ws.Range("H9").Value = [double]"2.871036741680655e-13"
$ws.Range("I9").Value = [double]"2.871036741680655e-13"
$ws.Range("L9").Value = [double]"58.51066702663084"
$ws.Range("M9").Value = "[44.233884771673175, 72.78744928158851]"
$ws.Range("N9").Value = [double]"1.473792199391255e-10"
$ws.Range("O9").Value = [double]"1.473792199391255e-10"
$ws.Range("P9").Value = [double]"-0.5660527303765397"
$ws.Range("Q9").Value = "[-0.8176317216550011, -0.3144737390980783]"
$ws.Range("R9").Value = [double]"4.294577757657514e-05"
$ws.Range("S9").Value = [double]"4.294577757657514e-05"
$ws.Range("T9").Value = [double]"50.99573008041484"
$ws.Range("U9").Value = "[43.24998987921005, 58.74147028161963]"
$ws.Range("V9").Value = [double]"0"
$ws.Range("W9").Value = [double]"0"
$ws.Range("X9").Value = [double]"2.322522522522579"
$ws.Range("Y9").Value = [double]"1.290290290290325"
$ws.Range("Z9").Value = [double]"3.354754754754833"

# Row 10
$ws.Range("F10").Value = [double]"25.78000000000059"
$ws.Range("H10").Value = [double]"3.130384840233091e-12"
$ws.Range("I10").Value = [double]"3.130384840233091e-12"
$ws.Range("L10").Value = [double]"55.07847438802241"
$ws.Range("M10").Value = "[39.83819773314855, 70.31875104289627]"
$ws.Range("N10").Value = [double]"3.937494774675088e-09"
$ws.Range("O10").Value = [double]"3.937494774675088e-09"
$ws.Range("P10").Value = [double]"-1.044052813805617"
$ws.Range("Q10").Value = "[-1.3207897042119257, -0.767315923399309]"
$ws.Range("R10").Value = [double]"1.331338816967786e-09"
$ws.Range("S10").Value = [double]"1.331338816967786e-09"
$ws.Range("T10").Value = [double]"50.9543976169638"
$ws.Range("U10").Value = "[42.86395983319857, 59.044835400729035]"
$ws.Range("V10").Value = [double]"2.220446049250313e-16"
$ws.Range("W10").Value = [double]"2.220446049250313e-16"
$ws.Range("X10").Value = [double]"4.283763763763865"
$ws.Range("Y10").Value = [double]"3.148308308308382"
$ws.Range("Z10").Value = [double]"5.419219219219348"
